# Update Wnt1-Fzd4 LR-pair sheet with new TPM-derived NATMI output.
#
# The new run only has two sending/target clusters (ECs, FAPs) crossed with
# three target clusters (ECs, FAPs, MuSCs) for the Wnt1 -> Fzd4 edge, i.e. a
# 2x3 = 6 row grid (rows 2-7). The previous 9-row sheet (3 senders x 3
# targets incl. "MuSCs" as a sender and "Resolving-Mac" as a target) is
# replaced: rows 8-9 are dropped and every remaining data row/column is
# rewritten with the refreshed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old sheet had 8 data rows (rows 2-9); the refreshed export only has 6
# (rows 2-7), so drop the trailing two rows entirely.
$ws.Range("A8:T9").EntireRow.Delete()

# New 2 (sender) x 3 (target) grid of rows, each followed by the 16 numeric
# metric columns E..T in order.
$rows = @(
    @{ Row=2; A="ECs";  D="ECs";   E=1; F=0.3333333333333333; G=0.131499;            H=0.394497;           I=0.3654391092296077; J=0.3654391092296077; K=3; L=1; M=34.97976933333334;  N=104.939308;         O=0.4352965780925344; P=0.4352965780925344; Q=4.599804687564001; R=41.398242188076;    S=0.1590743937488321; T=0.1590743937488321 },
    @{ Row=3; A="ECs";  D="FAPs";  E=1; F=0.3333333333333333; G=0.131499;            H=0.394497;           I=0.3654391092296077; J=0.3654391092296077; K=3; L=1; M=20.343383;           N=61.03014900000001;  O=0.2531579017099818; P=0.2531579017099818; Q=2.675134521117001; R=24.076210690053;    S=0.09251379809533233; T=0.09251379809533233 },
    @{ Row=4; A="ECs";  D="MuSCs"; E=1; F=0.3333333333333333; G=0.131499;            H=0.394497;           I=0.3654391092296077; J=0.3654391092296077; K=3; L=1; M=25.035323;           N=75.105969;          O=0.3115455201974837; P=0.3115455201974837; Q=3.292119939177;    R=29.629079452593;    S=0.1138509173854432; T=0.1138509173854432 },
    @{ Row=5; A="FAPs"; D="ECs";   E=1; F=0.3333333333333333; G=0.2283393333333333;  H=0.685018;           I=0.6345608907703922; J=0.6345608907703922; K=3; L=1; M=34.97976933333334;  N=104.939308;         O=0.4352965780925344; P=0.4352965780925344; Q=7.987257209727112; R=71.88531488754401;  S=0.2762221843437023; T=0.2762221843437022 },
    @{ Row=6; A="FAPs"; D="FAPs";  E=1; F=0.3333333333333333; G=0.2283393333333333;  H=0.685018;           I=0.6345608907703922; J=0.6345608907703922; K=3; L=1; M=20.343383;           N=61.03014900000001;  O=0.2531579017099818; P=0.2531579017099818; Q=4.645194511964667; R=41.80675060768201;  S=0.1606441036146495; T=0.1606441036146495 },
    @{ Row=7; A="FAPs"; D="MuSCs"; E=1; F=0.3333333333333333; G=0.2283393333333333;  H=0.685018;           I=0.6345608907703922; J=0.6345608907703922; K=3; L=1; M=25.035323;           N=75.105969;          O=0.3115455201974837; P=0.3115455201974837; Q=5.716548963604668; R=51.448940672442;    S=0.1976946028120405; T=0.1976946028120405 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = "Wnt1"
    $ws.Range("C$n").Value = "Fzd4"
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
    $ws.Range("I$n").Value = $r.I
    $ws.Range("J$n").Value = $r.J
    $ws.Range("K$n").Value = $r.K
    $ws.Range("L$n").Value = $r.L
    $ws.Range("M$n").Value = $r.M
    $ws.Range("N$n").Value = $r.N
    $ws.Range("O$n").Value = $r.O
    $ws.Range("P$n").Value = $r.P
    $ws.Range("Q$n").Value = $r.Q
    $ws.Range("R$n").Value = $r.R
    $ws.Range("S$n").Value = $r.S
    $ws.Range("T$n").Value = $r.T
}
